# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# sheets to reflect the latest scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value = 1350
$sheet1.Range("F8").Value = 1514
$sheet1.Range("F11").Value = 837
$sheet1.Range("F13").Value = 99
$sheet1.Range("F14").Value = 45
$sheet1.Range("F18").Value = 5977
$sheet1.Range("F20").Value = 5801
$sheet1.Range("F21").Value = 9790
$sheet1.Range("F25").Value = 265
$sheet1.Range("F29").Value = 4365

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F5").Value = 1350
$sheet4.Range("F12").Value = 1514
$sheet4.Range("F15").Value = 837
$sheet4.Range("F18").Value = 99
$sheet4.Range("F19").Value = 45
$sheet4.Range("F24").Value = 5977
$sheet4.Range("F26").Value = 5801
$sheet4.Range("F27").Value = 9790
$sheet4.Range("F32").Value = 265
$sheet4.Range("F39").Value = 4365
